$wb = $excel.ActiveWorkbook

# The workbook currently has two sheets: "Rutas" (1) and "Rutas_Unidad" (2).
# We need to insert a brand new sheet "Gastos por Unidad" with a CPK summary
# table as the new first sheet, pushing "Rutas" and "Rutas_Unidad" to
# positions 2 and 3 respectively (their own content stays unchanged).

$rutas = $wb.Worksheets.Item("Rutas")

$ws = $wb.Worksheets.Add($rutas)
$ws.Name = "Gastos por Unidad"

# Reuse the existing bold/bordered/centered header style from the "Rutas"
# sheet's header row (A1:D1) instead of building a brand-new style, so the
# underlying style table stays the same as before the edit.
$rutas.Range("A1:D1").Copy($ws.Range("A1:D1"))
$rutas.Range("A1:C1").Copy($ws.Range("E1:G1"))

$headers = @("Unidad", "Litros", "Gasto Combustible", "Costo por litro", "Kms Totales", "Gasto Mantenimiento", "CPK")
for ($c = 1; $c -le 7; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Data rows.
$data = @(
    @("1823", 23564859.44, 461196837.92868, 19.56694042857142, 1734208, 50879.31, 265.9702395783435),
    @("1839", 19890344.8, 389204913.07396, 19.57580044444443, 1534842, 131218.87, 253.665284077423),
    @("1885", 17639860.11, 346288330.362969, 19.61239606451614, 1717121, 10330.47, 201.6740001624632),
    @("1752", 20057740.62, 393427926.476439, 19.60499997999999, 2613450, 42203.16, 150.5558283634426),
    @("1970", 22033674.95, 436867450.244855, 20.1538869117647, 9793020, 39123.91, 44.61407963578702),
    @("1903", 19654148.73, 387318956.4651538, 20.07962663157895, 9165448, 14954.2, 42.2602267412519),
    @("TT02", 2754.76, 70508.1924, 25.83625000000001, 0, 5792.96, 0),
    @("DC05", 665.6900000000001, 12936.40969, 19.433084, 0, 99431.19, 0),
    @("DC04", 11298, 218873.21652, 19.37274, 0, 77044.37000000001, 0),
    @("DC02", 1239.52, 24087.69628, 19.433084, 0, 81901.66, 0)
)

$r = 2
foreach ($row in $data) {
    # Column A ("Unidad") holds unit codes that can look numeric (e.g.
    # "1823"), so force them to remain plain text the same way Excel's
    # UI would (leading apostrophe), matching the default (un-styled)
    # text cells used elsewhere in the workbook. Resetting the style to
    # "Normal" afterwards drops the quote-prefix style Excel applies,
    # so no new entries are added to the style table.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 1).Style = "Normal"
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r = $r + 1
}

Write-Host "Added sheet 'Gastos por Unidad' with" ($r - 2) "data rows."
